# Replace the hard-coded "Columbia, Missouri" location text with a
# {location} placeholder, in each of the three (Work / Volunteer / Extra)
# entries. The original run held both the <w:tab/> and the
# "Columbia, Missouri" text:
#
#   <w:r><w:rPr>...Georgia...</w:rPr><w:tab/><w:t>Columbia, Missouri</w:t></w:r>
#
# The target shape keeps the <w:tab/> in its own run and puts the new
# {location} placeholder text in a second, sibling run with the same
# run properties:
#
#   <w:r><w:rPr>...Georgia...</w:rPr><w:tab/></w:r>
#   <w:r><w:rPr>...Georgia...</w:rPr><w:t>{location}</w:t></w:r>
#
# InsertXML (unlike setting .Text) lets us drop exact OOXML into the
# document, so we delete the "<tab>Columbia, Missouri" span and replace it
# with those two runs verbatim.

$d = $word.ActiveDocument

$replacementXml = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
  <pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
    <pkg:xmlData>
      <w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
        <w:body>
          <w:p>
            <w:r>
              <w:rPr>
                <w:rFonts w:ascii="Georgia" w:hAnsi="Georgia"/>
              </w:rPr>
              <w:tab/>
            </w:r>
            <w:r>
              <w:rPr>
                <w:rFonts w:ascii="Georgia" w:hAnsi="Georgia"/>
              </w:rPr>
              <w:t>{location}</w:t>
            </w:r>
          </w:p>
        </w:body>
      </w:document>
    </pkg:xmlData>
  </pkg:part>
</pkg:package>
'@

# Repeatedly locate a paragraph still containing the literal text and fix
# it up; re-scanning after every edit keeps paragraph/range offsets valid
# since InsertXML changes the document length.
$keepGoing = $true
while ($keepGoing) {
    $keepGoing = $false
    foreach ($p in $d.Paragraphs) {
        $paraText = $p.Range.Text
        if ($paraText -like "*Columbia, Missouri*") {
            $paraStart = $p.Range.Start
            $paraEnd = $p.Range.End
            $tabIndex = $paraText.IndexOf("`t")

            # Span from the tab character through to just before the
            # paragraph mark; this is exactly the run that currently holds
            # "<w:tab/><w:t>Columbia, Missouri</w:t>".
            $target = $d.Range($paraStart + $tabIndex, $paraEnd - 1)
            $target.InsertXML($replacementXml)

            $keepGoing = $true
            break
        }
    }
}
